# Auto-generated edit script applying the Phantom_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N)
# across rows on sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 748.5
$ws.Range("I31").Value = 141.14285
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 423.42855
$ws.Range("L31").Value = 15000
$ws.Range("M31").Value = -193.42855
$ws.Range("N31").Value = -15460

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4919.4
$ws.Range("J76").Value = 3800
$ws.Range("L76").Value = 3800
$ws.Range("N76").Value = -4430

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4919.4
$ws.Range("J79").Value = 3800
$ws.Range("L79").Value = 3800
$ws.Range("N79").Value = -5984

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1143.2
$ws.Range("J92").Value = 1220
$ws.Range("L92").Value = 1220
$ws.Range("N92").Value = -3716

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 15152933
$ws.Range("I137").Value = 23810596
$ws.Range("J137").Value = 2024.625
$ws.Range("K137").Value = 71431788
$ws.Range("L137").Value = 6073.875
$ws.Range("M137").Value = -71429238
$ws.Range("N137").Value = -11173.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 5626.6665
$ws.Range("I31").Value = 5626.6665
$ws.Range("K31").Value = 5626.6665
$ws.Range("M31").Value = -5332.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6447.2354
$ws.Range("I32").Value = 4288.467
$ws.Range("K32").Value = 4288.467
$ws.Range("M32").Value = -4001.467

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 999
$ws.Range("I45").Value = 999
$ws.Range("K45").Value = 999
$ws.Range("M45").Value = -622

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4999
$ws.Range("I61").Value = 4999
$ws.Range("K61").Value = 4999
$ws.Range("M61").Value = -4787

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 813.8570999999999
$ws.Range("I97").Value = 763.9474
$ws.Range("J97").Value = 1288
$ws.Range("K97").Value = 763.9474
$ws.Range("L97").Value = 1288
$ws.Range("M97").Value = -267.9474
$ws.Range("N97").Value = -2280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4999
$ws.Range("I136").Value = 4999
$ws.Range("K136").Value = 14997
$ws.Range("M136").Value = -12447

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 12365487
$ws.Range("I86").Value = 20480.467
$ws.Range("J86").Value = 27796746
$ws.Range("K86").Value = 20480.467
$ws.Range("L86").Value = 27796746
$ws.Range("M86").Value = -19357.467
$ws.Range("N86").Value = -27798992

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 12776.8
$ws.Range("J88").Value = 12776.8
$ws.Range("L88").Value = 12776.8
$ws.Range("N88").Value = -13588.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 12365487
$ws.Range("I89").Value = 20480.467
$ws.Range("J89").Value = 27796746
$ws.Range("K89").Value = 102402.335
$ws.Range("L89").Value = 138983730
$ws.Range("M89").Value = -96786.33500000001
$ws.Range("N89").Value = -138994962

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H91").Value = 12776.8
$ws.Range("J91").Value = 12776.8
$ws.Range("L91").Value = 12776.8
$ws.Range("N91").Value = -15584.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3448.5
$ws.Range("I99").Value = 3448.5
$ws.Range("K99").Value = 3448.5
$ws.Range("M99").Value = -1950.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1994.5
$ws.Range("I134").Value = 1994.6666
$ws.Range("J134").Value = 1994
$ws.Range("K134").Value = 5983.9998
$ws.Range("L134").Value = 5982
$ws.Range("M134").Value = -3448.9998
$ws.Range("N134").Value = -11052

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 240.4
$ws.Range("I16").Value = 264.5
$ws.Range("K16").Value = 264.5
$ws.Range("M16").Value = 22.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 2000
$ws.Range("J29").Value = 2000
$ws.Range("L29").Value = 2000
$ws.Range("N29").Value = -2586

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2481
$ws.Range("I58").Value = 2920.3125
$ws.Range("K58").Value = 2920.3125
$ws.Range("M58").Value = -2717.3125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 50000
$ws.Range("I103").Value = 50000
$ws.Range("K103").Value = 50000
$ws.Range("M103").Value = -48828

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 240.4
$ws.Range("I113").Value = 264.5
$ws.Range("K113").Value = 264.5
$ws.Range("M113").Value = 1905.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2900
$ws.Range("I122").Value = 2900
$ws.Range("K122").Value = 8700
$ws.Range("M122").Value = -6250

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2481
$ws.Range("I136").Value = 2920.3125
$ws.Range("K136").Value = 8760.9375
$ws.Range("M136").Value = -6210.9375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 336410.2
$ws.Range("J141").Value = 336410.2
$ws.Range("L141").Value = 336410.2
$ws.Range("N141").Value = -346770.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 1003.2857
$ws.Range("I10").Value = 46
$ws.Range("K10").Value = 138
$ws.Range("M10").Value = 1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1580.5555
$ws.Range("I11").Value = 1841
$ws.Range("K11").Value = 5523
$ws.Range("M11").Value = -5383

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J25").Value = 500
$ws.Range("L25").Value = 1500
$ws.Range("N25").Value = -1838

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J30").Value = 500
$ws.Range("L30").Value = 1500
$ws.Range("N30").Value = -1704

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 10608.728
$ws.Range("I94").Value = 566.3333
$ws.Range("K94").Value = 1698.9999
$ws.Range("M94").Value = -1022.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 2151.3044
$ws.Range("I99").Value = 749.25
$ws.Range("K99").Value = 2247.75
$ws.Range("M99").Value = -1.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 1000
$ws.Range("I102").Value = 1000
$ws.Range("K102").Value = 3000
$ws.Range("M102").Value = -566

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 3333
$ws.Range("J127").Value = 3333
$ws.Range("L127").Value = 9999
$ws.Range("N127").Value = -19919

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 772465.75
$ws.Range("I140").Value = 772465.75
$ws.Range("K140").Value = 2317397.25
$ws.Range("M140").Value = -2312217.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3125.75
$ws.Range("I80").Value = 2499
$ws.Range("J80").Value = 5006
$ws.Range("K80").Value = 2499
$ws.Range("L80").Value = 5006
$ws.Range("M80").Value = -1501
$ws.Range("N80").Value = -7002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3125.75
$ws.Range("I83").Value = 2499
$ws.Range("J83").Value = 5006
$ws.Range("K83").Value = 12495
$ws.Range("L83").Value = 25030
$ws.Range("M83").Value = -7503
$ws.Range("N83").Value = -35014

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1645.2142
$ws.Range("I97").Value = 1645.2142
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1645.2142
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -1149.2142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2179.2
$ws.Range("I102").Value = 2186.6667
$ws.Range("K102").Value = 2186.6667
$ws.Range("M102").Value = -564.6667000000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3318.1
$ws.Range("I113").Value = 2297
$ws.Range("K113").Value = 2297
$ws.Range("M113").Value = -127

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 25643480
$ws.Range("I132").Value = 2676.889
$ws.Range("K132").Value = 8030.667
$ws.Range("M132").Value = -5500.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2999
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 528.3182
$ws.Range("I16").Value = 539.1905
$ws.Range("J16").Value = 300
$ws.Range("K16").Value = 539.1905
$ws.Range("L16").Value = 300
$ws.Range("M16").Value = -369.1905
$ws.Range("N16").Value = -640

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 71435690
$ws.Range("I40").Value = 166671790
$ws.Range("J40").Value = 8612.25
$ws.Range("K40").Value = 166671790
$ws.Range("L40").Value = 8612.25
$ws.Range("M40").Value = -166671654
$ws.Range("N40").Value = -8884.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H57").Value = 25000
$ws.Range("I57").Value = 25000
$ws.Range("K57").Value = 25000
$ws.Range("M57").Value = -24434

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3379.7856
$ws.Range("I61").Value = 3536.6667
$ws.Range("K61").Value = 3536.6667
$ws.Range("M61").Value = -3334.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1147.6875
$ws.Range("I82").Value = 847.4286
$ws.Range("K82").Value = 847.4286
$ws.Range("M82").Value = -486.4286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1147.6875
$ws.Range("I85").Value = 847.4286
$ws.Range("K85").Value = 847.4286
$ws.Range("M85").Value = 400.5714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3379.7856
$ws.Range("I113").Value = 3536.6667
$ws.Range("K113").Value = 3536.6667
$ws.Range("M113").Value = -1366.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2093.75
$ws.Range("I122").Value = 2093.75
$ws.Range("K122").Value = 6281.25
$ws.Range("M122").Value = -3831.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2999
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 13596.8
$ws.Range("I132").Value = 4496
$ws.Range("K132").Value = 13488
$ws.Range("M132").Value = -10958

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4248.875
$ws.Range("I136").Value = 4248.875
$ws.Range("K136").Value = 12746.625
$ws.Range("M136").Value = -10196.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 90910730
$ws.Range("I132").Value = 1778.3334
$ws.Range("J132").Value = 500001000
$ws.Range("K132").Value = 5335.0002
$ws.Range("L132").Value = 1500003000
$ws.Range("M132").Value = -2805.0002
$ws.Range("N132").Value = -1500008060
